# Aggiunti indicatori BIB livello 3
# Adds 11 new "Library_Formula" rows (CREATE/MODIFY | LIB_EWS_BE | INDICATOR_xxx | <blank D> | String)
# right after the existing data (rows 2-99), for the new level-3 BIB indicators.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

$indicators = @(
    "INDICATOR_8",
    "INDICATOR_13",
    "INDICATOR_59",
    "INDICATOR_61",
    "INDICATOR_62",
    "INDICATOR_63",
    "INDICATOR_67",
    "INDICATOR_68",
    "INDICATOR_172",
    "INDICATOR_187",
    "INDICATOR_188"
)

$row = 100
foreach ($ind in $indicators) {
    $cA = $ws.Cells.Item($row, 1)
    $cA.Value = "CREATE/MODIFY"
    $cA.Font.Name = "Trebuchet MS"
    $cA.Font.Size = 10

    $cB = $ws.Cells.Item($row, 2)
    $cB.Value = "LIB_EWS_BE"
    $cB.Font.Name = "Trebuchet MS"
    $cB.Font.Size = 10

    $cC = $ws.Cells.Item($row, 3)
    $cC.Value = $ind
    $cC.Font.Name = "Trebuchet MS"
    $cC.Font.Size = 10

    $cE = $ws.Cells.Item($row, 5)
    $cE.Value = "String"
    $cE.Font.Name = "Trebuchet MS"
    $cE.Font.Size = 10

    $row = $row + 1
}

# Reflect the author's final cursor position/selection on the sheet.
$ws.Range("C112").Select() | Out-Null
